# Update 20180925 - optiqueDesTissus/dataOtptiqueDesTissus.xlsx
# Adds a "semaine 2" (avec/sans correction) results table below the
# existing data, in rows 36-39, columns B-L.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- widen / add the columns used by the new table -------------------
$ws.Columns.Item(7).ColumnWidth  = 23.333333333333336
$ws.Columns.Item(10).ColumnWidth = 21.5
$ws.Columns.Item(11).ColumnWidth = 22.166666666666668
$ws.Columns.Item(12).ColumnWidth = 27.5

# --- row 36: section titles ------------------------------------------
# Text entered in this order so new shared-string entries land in the
# same order as the target workbook (semaine2, avec correction, sans
# correction, ...).
$ws.Range("B36").Value = "semaine 2"
$ws.Range("I36").Value = "avec correction"
$ws.Range("D36").Value = "sans correction"
$ws.Range("D36").NumberFormat = $ws.Range("G2").NumberFormat

# --- row 37: column headers for the two sub-tables --------------------
$ws.Range("C37").Value = "transmitance"
$ws.Range("D37").Value = "Réflectance"
$ws.Range("H37").Value = "tramsitance"
$ws.Range("I37").Value = "réflectance"
$ws.Range("E37").Value = "absorption (mm^-1)"
$ws.Range("F37").Value = "diffusion (mm^-1)"
$ws.Range("G37").Value = "diffusion reduce (mm^-1)"
$ws.Range("J37").Value = "absorption (mm^-1)"
$ws.Range("K37").Value = "diffusion (mm^-1)"
$ws.Range("L37").Value = "diffusion reduce (mm^-1)"
$ws.Range("G37").NumberFormat = $ws.Range("G2").NumberFormat
$ws.Range("I37").NumberFormat = $ws.Range("G2").NumberFormat
$ws.Range("L37").NumberFormat = $ws.Range("G2").NumberFormat

# --- row 38: raw measurements ------------------------------------------
$ws.Range("C38").Value = 0.2864
$ws.Range("D38").Value = 0.1898
$ws.Range("H38").Value = 0.3004
$ws.Range("I38").Value = 0.1932
$ws.Range("C38").NumberFormat = $ws.Range("G2").NumberFormat
$ws.Range("D38").NumberFormat = $ws.Range("G2").NumberFormat

# --- row 39: computed / entered results --------------------------------
$ws.Range("C39").Formula = "=C38/(C38+D38)"
$ws.Range("D39").Formula = "=D38/(D38/C38)"
$ws.Range("E39").Value = 0.06097
$ws.Range("F39").Value = 9.6843
$ws.Range("G39").Value = 1.227
$ws.Range("H39").Formula = "=H38/(H38+I38)"
$ws.Range("I39").Formula = "=I38/(H38+I38)"
$ws.Range("J39").Value = 0.00001246
$ws.Range("K39").Value = 9.7457
$ws.Range("L39").Value = 1.801

$ws.Range("C39,D39,E39,F39,G39,J39,K39,L39").NumberFormat = $ws.Range("G2").NumberFormat

# --- selection / active cell -------------------------------------------
$ws.Range("K42").Select()
